{"js": "// Update the title date and all 100 arithmetic-expression table cells to\n// match the \"within100\" worksheet's next edition.\n//\n// `replacements` lists, in document order, [originalText, newText] pairs:\n//   - index 0 is the centered title paragraph\n//     (\"2023-12-03 Sunday\" -> \"2023-12-04 Monday\")\n//   - indices 1..100 are the 20-row x 5-column table of \"a+b=\" / \"a-b=\"\n//     expressions, read left-to-right then top-to-bottom (same order as\n//     `body.paragraphs`, which walks the document including table cells).\nconst replacements = [\n  [\"2023-12-03 Sunday\", \"2023-12-04 Monday\"],\n  [\"84+7=\", \"69-29=\"],\n  [\"73-70=\", \"20+2=\"],\n  [\"90-84=\", \"57-25=\"],\n  [\"84-20=\", \"8+91=\"],\n  [\"58-44=\", \"73-1=\"],\n  [\"52-29=\", \"3+24=\"],\n  [\"57-50=\", \"88-20=\"],\n  [\"14+13=\", \"90-30=\"],\n  [\"40+1=\", \"25+64=\"],\n  [\"17+59=\", \"11+78=\"],\n  [\"26+60=\", \"64+1=\"],\n  [\"33-14=\", \"5+50=\"],\n  [\"26+33=\", \"16+25=\"],\n  [\"88-60=\", \"20+6=\"],\n  [\"61+5=\", \"2+86=\"],\n  [\"79-24=\", \"33+38=\"],\n  [\"71-39=\", \"89-79=\"],\n  [\"34+54=\", \"59+32=\"],\n  [\"23-20=\", \"28+3=\"],\n  [\"59+4=\", \"83+3=\"],\n  [\"59-39=\", \"35-1=\"],\n  [\"76-13=\", \"71-33=\"],\n  [\"40+0=\", \"94-52=\"],\n  [\"68-45=\", \"75-72=\"],\n  [\"7+26=\", \"78-71=\"],\n  [\"1+44=\", \"18+32=\"],\n  [\"72-31=\", \"25+66=\"],\n  [\"49+17=\", \"9-6=\"],\n  [\"81-31=\", \"11-1=\"],\n  [\"80+15=\", \"56-2=\"],\n  [\"82+3=\", \"34+30=\"],\n  [\"75-64=\", \"15+4=\"],\n  [\"45-13=\", \"13+57=\"],\n  [\"99-81=\", \"25+42=\"],\n  [\"57-21=\", \"7+77=\"],\n  [\"19-19=\", \"82+15=\"],\n  [\"56+1=\", \"88-20=\"],\n  [\"3+32=\", \"8+18=\"],\n  [\"69-2=\", \"23-18=\"],\n  [\"44+18=\", \"25-9=\"],\n  [\"50-44=\", \"14+44=\"],\n  [\"8+11=\", \"99-71=\"],\n  [\"85-59=\", \"31+63=\"],\n  [\"69-3=\", \"26-20=\"],\n  [\"98-11=\", \"7+75=\"],\n  [\"99-19=\", \"7+31=\"],\n  [\"66+30=\", \"48-10=\"],\n  [\"74-26=\", \"26+9=\"],\n  [\"18+40=\", \"91-75=\"],\n  [\"9+44=\", \"50+48=\"],\n  [\"54-25=\", \"37-30=\"],\n  [\"5+80=\", \"16+72=\"],\n  [\"15+8=\", \"52+11=\"],\n  [\"33-25=\", \"72-68=\"],\n  [\"23-11=\", \"52+35=\"],\n  [\"20+74=\", \"73+1=\"],\n  [\"39-31=\", \"96-44=\"],\n  [\"47-14=\", \"18+11=\"],\n  [\"6+45=\", \"91-58=\"],\n  [\"72+10=\", \"69-62=\"],\n  [\"36+35=\", \"9+9=\"],\n  [\"58-6=\", \"74+21=\"],\n  [\"17+30=\", \"78+3=\"],\n  [\"54-33=\", \"18+8=\"],\n  [\"62-29=\", \"98+0=\"],\n  [\"19+46=\", \"45-26=\"],\n  [\"74+6=\", \"86-58=\"],\n  [\"62-18=\", \"40+2=\"],\n  [\"95-62=\", \"97-27=\"],\n  [\"2+49=\", \"93-50=\"],\n  [\"94-63=\", \"56-19=\"],\n  [\"45+12=\", \"58-16=\"],\n  [\"41+17=\", \"20+68=\"],\n  [\"33-17=\", \"5+60=\"],\n  [\"44-1=\", \"16+7=\"],\n  [\"28+1=\", \"31-26=\"],\n  [\"34+38=\", \"16+35=\"],\n  [\"9+38=\", \"41+56=\"],\n  [\"83-82=\", \"38+0=\"],\n  [\"34+13=\", \"60-52=\"],\n  [\"49+29=\", \"22-15=\"],\n  [\"97-70=\", \"38-38=\"],\n  [\"23+9=\", \"82-72=\"],\n  [\"61+16=\", \"8+31=\"],\n  [\"99-17=\", \"72-6=\"],\n  [\"15+72=\", \"16+33=\"],\n  [\"47-21=\", \"23+61=\"],\n  [\"87-80=\", \"7-1=\"],\n  [\"23+46=\", \"99-58=\"],\n  [\"68-64=\", \"72+25=\"],\n  [\"8+25=\", \"52-50=\"],\n  [\"26+48=\", \"93+2=\"],\n  [\"38-22=\", \"63-15=\"],\n  [\"59-14=\", \"39-23=\"],\n  [\"48-2=\", \"46+19=\"],\n  [\"41+47=\", \"99-22=\"],\n  [\"5+20=\", \"67+32=\"],\n  [\"46-23=\", \"55-18=\"],\n  [\"34+59=\", \"26+13=\"],\n  [\"89+7=\", \"69-64=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== replacements.length) {\n  throw new Error(\n    \"Unexpected paragraph count: found \" + items.length +\n    \", expected \" + replacements.length\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const paragraph = items[i];\n  // Replace the whole paragraph's text in place so the run's formatting\n  // (font, size, etc.) carries over to the new text.\n  paragraph.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the title date and all 100 arithmetic-expression table cells to\n# match the \"within100\" worksheet's next edition.\n#\n# $titleOld/$titleNew cover the centered heading paragraph\n# (\"2023-12-03 Sunday\" -> \"2023-12-04 Monday\"). $cellReplacements lists, in\n# row-major order (row 1 col 1..5, row 2 col 1..5, ...), each table cell's\n# current [oldText, newText] pair -- the same order Cell($r, $c) walks the\n# 20-row x 5-column table.\n$d = $word.ActiveDocument\n\n$titleOld = '2023-12-03 Sunday'\n$titleNew = '2023-12-04 Monday'\n\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.Text = $titleNew\n\n$cellReplacements = @(\n    @('84+7=', '69-29='),\n    @('73-70=', '20+2='),\n    @('90-84=', '57-25='),\n    @('84-20=', '8+91='),\n    @('58-44=', '73-1='),\n    @('52-29=', '3+24='),\n    @('57-50=', '88-20='),\n    @('14+13=', '90-30='),\n    @('40+1=', '25+64='),\n    @('17+59=', '11+78='),\n    @('26+60=', '64+1='),\n    @('33-14=', '5+50='),\n    @('26+33=', '16+25='),\n    @('88-60=', '20+6='),\n    @('61+5=', '2+86='),\n    @('79-24=', '33+38='),\n    @('71-39=', '89-79='),\n    @('34+54=', '59+32='),\n    @('23-20=', '28+3='),\n    @('59+4=', '83+3='),\n    @('59-39=', '35-1='),\n    @('76-13=', '71-33='),\n    @('40+0=', '94-52='),\n    @('68-45=', '75-72='),\n    @('7+26=', '78-71='),\n    @('1+44=', '18+32='),\n    @('72-31=', '25+66='),\n    @('49+17=', '9-6='),\n    @('81-31=', '11-1='),\n    @('80+15=', '56-2='),\n    @('82+3=', '34+30='),\n    @('75-64=', '15+4='),\n    @('45-13=', '13+57='),\n    @('99-81=', '25+42='),\n    @('57-21=', '7+77='),\n    @('19-19=', '82+15='),\n    @('56+1=', '88-20='),\n    @('3+32=', '8+18='),\n    @('69-2=', '23-18='),\n    @('44+18=', '25-9='),\n    @('50-44=', '14+44='),\n    @('8+11=', '99-71='),\n    @('85-59=', '31+63='),\n    @('69-3=', '26-20='),\n    @('98-11=', '7+75='),\n    @('99-19=', '7+31='),\n    @('66+30=', '48-10='),\n    @('74-26=', '26+9='),\n    @('18+40=', '91-75='),\n    @('9+44=', '50+48='),\n    @('54-25=', '37-30='),\n    @('5+80=', '16+72='),\n    @('15+8=', '52+11='),\n    @('33-25=', '72-68='),\n    @('23-11=', '52+35='),\n    @('20+74=', '73+1='),\n    @('39-31=', '96-44='),\n    @('47-14=', '18+11='),\n    @('6+45=', '91-58='),\n    @('72+10=', '69-62='),\n    @('36+35=', '9+9='),\n    @('58-6=', '74+21='),\n    @('17+30=', '78+3='),\n    @('54-33=', '18+8='),\n    @('62-29=', '98+0='),\n    @('19+46=', '45-26='),\n    @('74+6=', '86-58='),\n    @('62-18=', '40+2='),\n    @('95-62=', '97-27='),\n    @('2+49=', '93-50='),\n    @('94-63=', '56-19='),\n    @('45+12=', '58-16='),\n    @('41+17=', '20+68='),\n    @('33-17=', '5+60='),\n    @('44-1=', '16+7='),\n    @('28+1=', '31-26='),\n    @('34+38=', '16+35='),\n    @('9+38=', '41+56='),\n    @('83-82=', '38+0='),\n    @('34+13=', '60-52='),\n    @('49+29=', '22-15='),\n    @('97-70=', '38-38='),\n    @('23+9=', '82-72='),\n    @('61+16=', '8+31='),\n    @('99-17=', '72-6='),\n    @('15+72=', '16+33='),\n    @('47-21=', '23+61='),\n    @('87-80=', '7-1='),\n    @('23+46=', '99-58='),\n    @('68-64=', '72+25='),\n    @('8+25=', '52-50='),\n    @('26+48=', '93+2='),\n    @('38-22=', '63-15='),\n    @('59-14=', '39-23='),\n    @('48-2=', '46+19='),\n    @('41+47=', '99-22='),\n    @('5+20=', '67+32='),\n    @('46-23=', '55-18='),\n    @('34+59=', '26+13='),\n    @('89+7=', '69-64=')\n)\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif ($rowCount * $colCount -ne $cellReplacements.Count) {\n    throw \"Unexpected table size: found $rowCount x $colCount cells, expected $($cellReplacements.Count) replacements.\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $pair = $cellReplacements[$i]\n        $newText = $pair[1]\n        $cell = $table.Cell($r, $c)\n        # Assigning Range.Text replaces the cell's content while keeping the\n        # run/paragraph formatting (font, size, etc.) already in the cell.\n        $cell.Range.Text = $newText\n        $i = $i + 1\n    }\n}\n\nWrite-Output \"Updated $i table cells and the title paragraph.\"\n"}
